# Word COM-interop script: update the "Please note" wording in the
# Notice of Trial hearing documents.
#
# The paragraph that used to read:
#   "Please note: This case may be released to a different court hearing
#   centre, in which case you will be notified."
# becomes:
#   "Please note: Cases are listed in accordance with local hearing
#   arrangements determined by the Judiciary and implemented by the court
#   staff. Every effort is made to ensure that hearings start at the time
#   specified. However, listing practices or other factors may mean that
#   you experience a delay, an adjournment at short notice or your case
#   may be released to a different court hearing centre, in which case
#   you will be notified."
#
# The bold "Please note: " lead-in run is left untouched; only the
# non-bold sentence(s) that follow it are replaced.

$d = $word.ActiveDocument

$oldText = "This case may be released to a different court hearing centre, in which case you will be notified."
$newText = "Cases are listed in accordance with local hearing arrangements determined by the Judiciary and implemented by the court staff. Every effort is made to ensure that hearings start at the time specified. However, listing practices or other factors may mean that you experience a delay, an adjournment at short notice or your case may be released to a different court hearing centre, in which case you will be notified."

$rng = $d.Content
$rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
